# Update "Förändrad" (column C) date values for rows 2-9 from 45208 to 45212
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45212
}
